$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before C. This shifts the old "Trait" column (C) to D,
#    and the old "Desc" column (D) to E.
$ws.Columns("C:C").Insert()

# 2. Clear the leftover formatting on the new column header (C1, copied from column B)
#    and on the shifted TraitName column (now D) so they return to the default
#    (unstyled) cell format before we re-apply the header's bold style below.
$ws.Range("C1").ClearFormats()
$ws.Range("D1:D25").ClearFormats()

# 3. Set the new header row text.
$ws.Range("C1").Value = "PersonalityType"
$ws.Range("D1").Value = "TraitName"

# 4. Make the two new/renamed header cells bold (Calibri 12, the workbook default font).
$ws.Range("C1:D1").Font.Bold = $true

# 5. Fill in the PersonalityType values for the 16Personalitties rows (2-17).
$ws.Range("C2").Value = "ENTJ"
$ws.Range("C3").Value = "EOTJ"
$ws.Range("C4").Value = "ENPJ"
$ws.Range("C5").Value = "EOPJ"
$ws.Range("C6").Value = "ENTF"
$ws.Range("C7").Value = "EOTF"
$ws.Range("C8").Value = "ENPF"
$ws.Range("C9").Value = "EOPF"
$ws.Range("C10").Value = "INTJ"
$ws.Range("C11").Value = "IOTJ"
$ws.Range("C12").Value = "INPJ"
$ws.Range("C13").Value = "IOPJ"
$ws.Range("C14").Value = "INTF"
$ws.Range("C15").Value = "IOTF"
$ws.Range("C16").Value = "INPF"
$ws.Range("C17").Value = "IOPF"

# 6. For the Love/Job rows (18-25), PersonalityType duplicates the TraitName value.
$ws.Range("C18").Value = "Quality Time"
$ws.Range("C19").Value = "Act of Service"
$ws.Range("C20").Value = "Physical Touch"
$ws.Range("C21").Value = "Words of Affirmation"
$ws.Range("C22").Value = "Outgoing"
$ws.Range("C23").Value = "Asocial"
$ws.Range("C24").Value = "Idealistic"
$ws.Range("C25").Value = "Realistic"

# 7. Match the column widths: PersonalityType (C) shares the TraitName (D) width.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# 8. Restore the active selection.
$ws.Range("D24").Select()
